# Applies the "Update pptx presentation file." commit:
#  - Slide 10 (Model evaluation): expand the three summary lines into
#    header + two sub-bullets (R-squared / RMSE) each, sub-bullets indented
#    one level.
#  - Slide 6 (Metrics and Limitations): extend the "Look at R-squared test
#    set values." bullet to also mention RMSE.

$p = $ppt.ActivePresentation
$cr = [char]13

# --- Slide 10: "Model evaluation" -----------------------------------------
$s10 = $p.Slides.Item(10)
$content10 = $s10.Shapes.Item(2)
$tr10 = $content10.TextFrame.TextRange

$tr10.Text = "Linear regression:  " + $cr + `
             "R-squared test value: 0.724" + $cr + `
             "RMSE test value: 2.28" + $cr + `
             "Decision tree:" + $cr + `
             "R-squared test value: 0.553" + $cr + `
             "RMSE test value: 2.90" + $cr + `
             "Random forest:" + $cr + `
             "R-squared test value: 0.700" + $cr + `
             "RMSE test value: 2.38"

# Indent the R-squared / RMSE detail paragraphs one level (lvl="1").
$tr10.Paragraphs(2,1).IndentLevel = 2
$tr10.Paragraphs(3,1).IndentLevel = 2
$tr10.Paragraphs(5,1).IndentLevel = 2
$tr10.Paragraphs(6,1).IndentLevel = 2
$tr10.Paragraphs(8,1).IndentLevel = 2
$tr10.Paragraphs(9,1).IndentLevel = 2

# --- Slide 6: "Metrics and Limitations" ------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange

$rsquaredPara = $tr6.Paragraphs(2,1)
# Two-step set avoids the common prefix/suffix "." being kept as its own
# run - we want a single clean <a:r> for the whole new sentence.
$rsquaredPara.Text = "placeholder"
$tr6.Paragraphs(2,1).Text = "Look at R-squared test set values and RMSE values."
